$p = $ppt.ActivePresentation

# The deck currently has 3 slides: [logo] [card] [card].
# We duplicate the last slide (the "card"), move the duplicate so it
# becomes the new 3rd slide, and restyle it (background color + the two
# emoji text boxes repositioned to match the 2nd slide's layout). The
# untouched original slide is pushed down to become the new 4th slide.

$original = $p.Slides.Item(3)
$newSlide = $original.Duplicate()
$newSlide.MoveTo(3)

$newSlide.Background.Fill.ForeColor.RGB = 3329586  # &H32CE32 (RGB(0x32,0xCE,0x32))

$dollar = $newSlide.Shapes.Item(1)   # "TextBox 3" - the money-bag emoji
$dollar.Left = 422.59456692913386
$dollar.Top = 270.0
$dollar.Width = 254.59456692913386
$dollar.Height = 201.1452755905512

$baseball = $newSlide.Shapes.Item(2) # "TextBox 4" - the baseball emoji
$baseball.Left = 225.40547213336615
$baseball.Top = 118.70275590551181
$baseball.Width = 300.97303118848424
$baseball.Height = 201.1452755905512
